$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-11 Tuesday", "2024-06-12 Wednesday"),
    @("56×60=", "80×28="),
    @("98×99=", "60×61="),
    @("94×65=", "77×28="),
    @("90×75=", "96×26="),
    @("80×17=", "64×22="),
    @("35×47=", "72×99="),
    @("83×80=", "76×63="),
    @("49×17=", "78×46="),
    @("74×49=", "72×88="),
    @("61×66=", "41×76="),
    @("26×62=", "38×53="),
    @("47×94=", "34×13="),
    @("84×90=", "46×56="),
    @("75×33=", "55×53="),
    @("48×22=", "20×62="),
    @("66×46=", "65×34="),
    @("46×84=", "94×97="),
    @("93×42=", "73×29="),
    @("31×96=", "56×49="),
    @("22×75=", "69×34="),
    @("93×53=", "26×25="),
    @("48×61=", "29×61="),
    @("63×15=", "86×26="),
    @("19×73=", "54×98="),
    @("48×17=", "44×94=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
